$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Duplicate the "LargeBin" sheet (while it still holds its original,
#    unedited content) to create the new "StorageBox" sheet, placed
#    right after "LargeBin" and before "Keras Seq".
# ------------------------------------------------------------------
$largeBin = $wb.Worksheets.Item("LargeBin")
$largeBin.Copy($null, $largeBin)
$storageBox = $wb.Worksheets.Item(4)
$storageBox.Name = "StorageBox"

# Label the new sheet "Storage Box" (new shared string, will become
# the first newly-added string -> index 38).
$storageBox.Range("B4").Value = "Storage Box"

# Match the new sheet's selection state (no tabSelected, single merged
# selection on the key/legend header row).
$storageBox.Activate()
$storageBox.Range("A5:C5").Select()

# ------------------------------------------------------------------
# 2. Finish editing the original "LargeBin" sheet: rename its label,
#    fill in the two missing result rows, and update its selection.
# ------------------------------------------------------------------
$largeBin.Range("B4").Value = "Large Bin"

$largeBin.Range("I45").Value = 0.884
$largeBin.Range("J45").Value = 0.34
$largeBin.Range("I46").Value = 0.679
$largeBin.Range("J46").Value = 0.622

$largeBin.Activate()
$largeBin.Range("A5:C5").Select()

# ------------------------------------------------------------------
# 3. Move the active tab / selection back to "ClosedDoor" (it becomes
#    the selected tab when the workbook is reopened) while keeping its
#    own selection unchanged.
# ------------------------------------------------------------------
$closedDoor = $wb.Worksheets.Item("ClosedDoor")
$closedDoor.Activate()
$closedDoor.Range("J48").Select()
